$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# The (empty) column Q is removed; everything from column R onward
# (R:AI) shifts one column to the left (becomes Q:AH).
$ws.Range("Q:Q").Delete()

# Re-point the two "containsText" conditional-formatting rules (previously
# applied across S/U/W/Y/AA/AC/AH) onto the shifted columns (R/T/V/X/Z/AB/AG)
# and fix up their search formulas to reference the new left-most column.
$cfColl = $ws.Range("S1").FormatConditions
$ruleUpdated = $cfColl.Item(1)
$ruleFilledIn = $cfColl.Item(2)
$ruleUpdated.ModifyAppliesToRange($ws.Range("R1:R1048576"))

$cfColl2 = $ws.Range("R1").FormatConditions
$ruleUpdated2 = $cfColl2.Item(1)
$ruleFilledIn2 = $cfColl2.Item(2)
$ruleUpdated2.Formula1 = 'NOT(ISERROR(SEARCH("updated",R1)))'
$ruleFilledIn2.Formula1 = 'NOT(ISERROR(SEARCH("filled in",R1)))'

# Match the saved view state: whole used range selected, scrolled near column L.
$ws.Activate() | Out-Null
$ws.Range("A1:AH7").Select() | Out-Null
